# Correction texte : "prédire" -> "détecter" (+ mise en forme des deux légendes)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "ZoneTexte 5" (legende gauche : "Modèle d'entrainement créé pour ... 3 cas") ---
$shape11 = $s.Shapes.Item(11)
$tr11 = $shape11.TextFrame.TextRange
$sub11 = $tr11.Characters(1, $tr11.Length)
$sub11.Text = "Modèle d’entrainement créé pour détecter 3 cas"
$tr11.Font.Size = 10
$tr11.Font.Bold = $true
# Taille figee par spAutoFit : on repositionne la hauteur exacte apres la MAJ du texte/police
$shape11.Height = 19.387519685039372

# --- Shape "ZoneTexte 22" (legende droite : "Utilisation du modèle INCEPTIONV3 ...") ---
$shape15 = $s.Shapes.Item(15)
$tr15 = $shape15.TextFrame.TextRange
$sub15 = $tr15.Characters(1, $tr15.Length)
$sub15.Text = "Utilisation du modèle INCEPTIONV3 pour classer 24 maladies de peau"
$tr15.Font.Size = 10
$tr15.Font.Bold = $true
# Idem : on fixe largeur/hauteur finales apres MAJ du texte/police
$shape15.Width = 395.7488582677165
$shape15.Height = 19.387519685039372
